$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.755.05"
$ws.Range("E2").Value = "  +2.80%  "
$ws.Range("D3").Value = "2.098.09"
$ws.Range("E3").Value = "  +2.92%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "228.18"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("E6").Value = "  +1.12%  "
$ws.Range("D7").Value = "60.44"
$ws.Range("E7").Value = "  +1.63%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +2.14%  "
$ws.Range("D10").Value = "0.0836"
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").Value = "2.409.86"
$ws.Range("E12").Value = "  +3.08%  "
$ws.Range("D13").Value = "14.97"
$ws.Range("E13").Value = "  +3.84%  "
$ws.Range("D14").Value = "22.22"
$ws.Range("E14").Value = "  +5.57%  "
$ws.Range("D15").Value = "0.793"
$ws.Range("E15").Value = "  +2.83%  "
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "2.098.04"
$ws.Range("E17").Value = "  +3.40%  "
$ws.Range("D18").Value = "38.735.94"
$ws.Range("E18").Value = "  +2.81%  "
$ws.Range("D19").Value = "71.66"
$ws.Range("E19").Value = "  +3.40%  "
$ws.Range("D20").Value = "6.04"
$ws.Range("E20").Value = "  +1.58%  "
$ws.Range("D21").Value = "0.0₃0834"
$ws.Range("E21").Value = "  +1.44%  "
$ws.Range("D22").Value = "225.69"
$ws.Range("E22").Value = "  +0.84%  "
$ws.Range("E24").Value = "  +1.49%  "
$ws.Range("E25").Value = "  +2.78%  "
$ws.Range("D26").Value = "170.74"
$ws.Range("E26").Value = "  +1.68%  "
$ws.Range("E27").Value = "  +1.10%  "
$ws.Range("D28").Value = "0.136"
$ws.Range("E28").Value = "  +5.80%  "
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").Value = "1.38"
$ws.Range("E29").Value = "  +8.84%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "19.18"
$ws.Range("E30").Value = "  +2.13%  "
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("E32").Value = "  +5.01%  "
$ws.Range("D33").Value = "4.76"
$ws.Range("E33").Value = "  +6.67%  "
$ws.Range("E34").Value = "  +2.92%  "
$ws.Range("D35").Value = "0.0611"
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("D36").Value = "2.39"
$ws.Range("E36").Value = "  +2.11%  "
$ws.Range("D37").Value = "6.41"
$ws.Range("E37").Value = "  -2.08%  "
$ws.Range("E38").Value = "  +3.50%  "
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("D40").Value = "18.39"
$ws.Range("E40").Value = "  +1.66%  "
$ws.Range("D41").Value = "1.543.96"
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("E42").Value = "  +4.35%  "
$ws.Range("D43").Value = "0.0222"
$ws.Range("E43").Value = "  +3.20%  "
$ws.Range("E44").Value = "  -0.76%  "
$ws.Range("D45").Value = "0.0925"
$ws.Range("E45").Value = "  +1.81%  "
$ws.Range("D46").Value = "7.66"
$ws.Range("E46").Value = "  +8.65%  "
$ws.Range("D47").Value = "4.13"
$ws.Range("E47").Value = "  -3.49%  "
$ws.Range("D48").Value = "1.11"
$ws.Range("E48").Value = "  +0.48%  "
$ws.Range("E49").Value = "  +2.49%  "
$ws.Range("D50").Value = "2.99"
$ws.Range("E50").Value = "  +1.68%  "
$ws.Range("D51").Value = "2.298.50"
$ws.Range("E51").Value = "  +3.14%  "
